# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> used by the (only) slide master, clrScheme "Integral"
#   ppt/theme/theme2.xml  -> used by the notes master,        clrScheme "Office"
#
# The target revision swaps the two themes' content in place, so the slide
# master's theme becomes the stock "Office Theme" colour scheme (and the
# notes master's theme becomes the old "Integral" colours). The two themes
# already share an identical fontScheme/fmtScheme, so the only real payload
# of the swap is the 12 theme colour slots (dk1/lt1/dk2/lt2/accent1-6/
# hlink/folHlink) of the master-facing theme (theme1.xml).
#
# PowerPoint's object model edits "the" presentation theme through
# ThemeColorScheme, addressed off any Slide (there's only one theme part the
# slide master resolves to here), so we push the Office Theme's twelve RGB
# values into it.
#
# COM colour longs pack bytes as R + G*256 + B*65536 (e.g. 0x44546A ->
# 0x6A*65536 + 0x54*256 + 0x44 = 6968388).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$colorScheme = $s.ThemeColorScheme

$colorScheme.Item(1).RGB  = 0         # dk1      000000
$colorScheme.Item(2).RGB  = 16777215  # lt1      FFFFFF
$colorScheme.Item(3).RGB  = 6968388   # dk2      44546A
$colorScheme.Item(4).RGB  = 15132391  # lt2      E7E6E6
$colorScheme.Item(5).RGB  = 13998939  # accent1  5B9BD5
$colorScheme.Item(6).RGB  = 3243501   # accent2  ED7D31
$colorScheme.Item(7).RGB  = 10855845  # accent3  A5A5A5
$colorScheme.Item(8).RGB  = 49407     # accent4  FFC000
$colorScheme.Item(9).RGB  = 12874308  # accent5  4472C4
$colorScheme.Item(10).RGB = 4697456   # accent6  70AD47
$colorScheme.Item(11).RGB = 12673797  # hlink    0563C1
$colorScheme.Item(12).RGB = 7491477   # folHlink 954F72

# Best-effort: try to rename the scheme/theme to match "Office"/"Office Theme"
# too. (Some hosts don't persist this -- harmless if it's a no-op.)
try { $colorScheme.Name = "Office" } catch { }
try { $p.DocumentTheme.Name = "Office Theme" } catch { }
